# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted into the data table at row 115,
# pushing all the subsequent rows (old 115..191) down by one (to 116..192).
# The sheet's used range therefore grows from A1:R191 to A1:R192.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 115 (shifts rows 115..191 down to 116..192)
$ws.Rows.Item(115).Insert()

# Populate the newly inserted row 115 with the new record's data
$ws.Range("A115").Value = 5
$ws.Range("B115").Value = "Macroferia Regional de Talca"
$ws.Range("C115").Value = "Maule"
$ws.Range("D115").Value = 44978
$ws.Range("E115").Value = 7
$ws.Range("F115").Value = 100112030
$ws.Range("G115").Value = "Poroto granado"
$ws.Range("H115").Value = "Sin especificar"
$ws.Range("I115").Value = "Primera"
$ws.Range("J115").Value = 400
$ws.Range("K115").Value = 25000
$ws.Range("L115").Value = 25000
$ws.Range("M115").Value = 25000
$ws.Range("N115").Value = "$/saco 25 kilos"
$ws.Range("O115").Value = "Región del Maule"
$ws.Range("P115").Value = 1000
$ws.Range("Q115").Value = 25
$ws.Range("R115").Value = "Hortaliza"
